# Auto-generated edit script applying numeric corrections to leve-profit calcs
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (scheduled price refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 754.6585
$ws.Range("I98").Value = 804.3158
$ws.Range("K98").Value = 804.3158
$ws.Range("M98").Value = 693.6842

$ws.Range("H106").Value = 8279.838
$ws.Range("I106").Value = 3015.4
$ws.Range("K106").Value = 3015.4
$ws.Range("M106").Value = -2384.4

$ws.Range("H111").Value = 841
$ws.Range("I111").Value = 711.625
$ws.Range("J111").Value = 1099.75
$ws.Range("K111").Value = 2134.875
$ws.Range("L111").Value = 3299.25
$ws.Range("M111").Value = 932.125
$ws.Range("N111").Value = -9433.25

$ws.Range("H116").Value = 35874.824
$ws.Range("I116").Value = 33738.363
$ws.Range("J116").Value = 37833.25
$ws.Range("K116").Value = 33738.363
$ws.Range("L116").Value = 37833.25
$ws.Range("M116").Value = -30296.363
$ws.Range("N116").Value = -44717.25

$ws.Range("H122").Value = 754.6585
$ws.Range("I122").Value = 804.3158
$ws.Range("K122").Value = 2412.9474
$ws.Range("M122").Value = 37.05259999999998

$ws.Range("H132").Value = 8006.794
$ws.Range("I132").Value = 8006.794
$ws.Range("K132").Value = 24020.382
$ws.Range("M132").Value = -21490.382

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2593.75
$ws.Range("I2").Value = 2458.3333
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 2458.3333
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -2345.3333
$ws.Range("N2").Value = -3226

$ws.Range("H5").Value = 1561.1538
$ws.Range("I5").Value = 177.33333
$ws.Range("J5").Value = 4674.75
$ws.Range("K5").Value = 177.33333
$ws.Range("L5").Value = 4674.75
$ws.Range("M5").Value = -65.33332999999999
$ws.Range("N5").Value = -4898.75

$ws.Range("H32").Value = 152777.7
$ws.Range("I32").Value = 180731.27
$ws.Range("K32").Value = 180731.27
$ws.Range("M32").Value = -180444.27

$ws.Range("H45").Value = 2458.8333
$ws.Range("I45").Value = 1788.25
$ws.Range("J45").Value = 3800
$ws.Range("K45").Value = 1788.25
$ws.Range("L45").Value = 3800
$ws.Range("M45").Value = -1411.25
$ws.Range("N45").Value = -4554

$ws.Range("H80").Value = 9900
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 9900
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H97").Value = 760.6061
$ws.Range("I97").Value = 712.9259
$ws.Range("J97").Value = 975.1667
$ws.Range("K97").Value = 712.9259
$ws.Range("L97").Value = 975.1667
$ws.Range("M97").Value = -216.9259
$ws.Range("N97").Value = -1967.1667

$ws.Range("H116").Value = 2593.75
$ws.Range("I116").Value = 2458.3333
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 2458.3333
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = -164.3332999999998
$ws.Range("N116").Value = -7588

$ws.Range("H132").Value = 611281.5600000001
$ws.Range("I132").Value = 626546.2
$ws.Range("K132").Value = 1879638.6
$ws.Range("M132").Value = -1877108.6

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2593.75
$ws.Range("I3").Value = 2458.3333
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 2458.3333
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -2344.3333
$ws.Range("N3").Value = -3228

$ws.Range("H4").Value = 1561.1538
$ws.Range("I4").Value = 177.33333
$ws.Range("J4").Value = 4674.75
$ws.Range("K4").Value = 177.33333
$ws.Range("L4").Value = 4674.75
$ws.Range("M4").Value = -62.33332999999999
$ws.Range("N4").Value = -4904.75

$ws.Range("H99").Value = 7402.3335

$ws.Range("H126").Value = 10000
$ws.Range("J126").Value = 10000
$ws.Range("L126").Value = 10000
$ws.Range("N126").Value = -19880

$ws.Range("H134").Value = 10777.714
$ws.Range("I134").Value = 3808
$ws.Range("K134").Value = 11424
$ws.Range("M134").Value = -8889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H25").Value = 4994
$ws.Range("I25").Value = 4994
$ws.Range("K25").Value = 4994
$ws.Range("M25").Value = -4820

$ws.Range("H107").Value = 2224
$ws.Range("I107").Value = 2224
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2224
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -304
$ws.Range("N107").ClearContents()

$ws.Range("H127").Value = 90000
$ws.Range("J127").Value = 90000
$ws.Range("L127").Value = 90000
$ws.Range("N127").Value = -99920

$ws.Range("H132").Value = 1610
$ws.Range("I132").Value = 1295.234
$ws.Range("K132").Value = 3885.702
$ws.Range("M132").Value = -1355.702

$ws.Range("H134").Value = 1143.8788
$ws.Range("I134").Value = 1153.1613
$ws.Range("K134").Value = 3459.4839
$ws.Range("M134").Value = -924.4839000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3612.2222
$ws.Range("I3").Value = 3438.75
$ws.Range("K3").Value = 10316.25
$ws.Range("M3").Value = -10204.25

$ws.Range("H15").Value = 50249.6
$ws.Range("I15").Value = 111400.78
$ws.Range("K15").Value = 334202.34
$ws.Range("M15").Value = -334062.34

$ws.Range("J34").Value = 5000
$ws.Range("L34").Value = 15000
$ws.Range("N34").Value = -15168

$ws.Range("H107").Value = 475.63635
$ws.Range("J107").Value = 1738
$ws.Range("L107").Value = 5214
$ws.Range("N107").Value = -9054

$ws.Range("H109").Value = 2713.889
$ws.Range("I109").Value = 1285
$ws.Range("J109").Value = 4500
$ws.Range("K109").Value = 3855
$ws.Range("L109").Value = 13500
$ws.Range("M109").Value = -2815
$ws.Range("N109").Value = -15580

$ws.Range("H137").Value = 4475.9
$ws.Range("J137").Value = 8666.666999999999
$ws.Range("L137").Value = 26000.001
$ws.Range("N137").Value = -36200.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2530.4285
$ws.Range("I43").Value = 2530.4285
$ws.Range("K43").Value = 2530.4285
$ws.Range("M43").Value = -2379.4285

$ws.Range("H80").Value = 1923.75
$ws.Range("I80").Value = 1798
$ws.Range("J80").Value = 2049.5
$ws.Range("K80").Value = 1798
$ws.Range("L80").Value = 2049.5
$ws.Range("M80").Value = -800
$ws.Range("N80").Value = -4045.5

$ws.Range("H82").Value = 30000
$ws.Range("I82").Value = 30000
$ws.Range("K82").Value = 30000
$ws.Range("M82").Value = -29617

$ws.Range("H83").Value = 1923.75
$ws.Range("I83").Value = 1798
$ws.Range("J83").Value = 2049.5
$ws.Range("K83").Value = 8990
$ws.Range("L83").Value = 10247.5
$ws.Range("M83").Value = -3998
$ws.Range("N83").Value = -20231.5

$ws.Range("H85").Value = 30000
$ws.Range("I85").Value = 30000
$ws.Range("K85").Value = 30000
$ws.Range("M85").Value = -28674

$ws.Range("H92").Value = 8658.333000000001
$ws.Range("J92").Value = 8658.333000000001
$ws.Range("L92").Value = 8658.333000000001
$ws.Range("N92").Value = -12402.333

$ws.Range("H113").Value = 2153.1052
$ws.Range("I113").Value = 2077.75
$ws.Range("J113").Value = 2555
$ws.Range("K113").Value = 2077.75
$ws.Range("L113").Value = 2555
$ws.Range("M113").Value = 92.25
$ws.Range("N113").Value = -6895

$ws.Range("H132").Value = 13735
$ws.Range("I132").Value = 14083.947
$ws.Range("K132").Value = 42251.841
$ws.Range("M132").Value = -39721.841

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2883.7273
$ws.Range("I22").Value = 1605.5555
$ws.Range("J22").Value = 3768.6155
$ws.Range("K22").Value = 1605.5555
$ws.Range("L22").Value = 3768.6155
$ws.Range("M22").Value = -1310.5555
$ws.Range("N22").Value = -4358.6155

$ws.Range("H27").Value = 2883.7273
$ws.Range("I27").Value = 1605.5555
$ws.Range("J27").Value = 3768.6155
$ws.Range("K27").Value = 1605.5555
$ws.Range("L27").Value = 3768.6155
$ws.Range("M27").Value = -1498.5555
$ws.Range("N27").Value = -3982.6155

$ws.Range("H82").Value = 2749.75
$ws.Range("I82").Value = 1999.6666
$ws.Range("K82").Value = 1999.6666
$ws.Range("M82").Value = -1638.6666

$ws.Range("H85").Value = 2749.75
$ws.Range("I85").Value = 1999.6666
$ws.Range("K85").Value = 1999.6666
$ws.Range("M85").Value = -751.6666

$ws.Range("H132").Value = 2647.1555
$ws.Range("I132").Value = 2352.718
$ws.Range("K132").Value = 7058.154
$ws.Range("M132").Value = -4528.154

$ws.Range("H136").Value = 1251.9778
$ws.Range("I136").Value = 1066.0233
$ws.Range("K136").Value = 3198.0699
$ws.Range("M136").Value = -648.0699000000004

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H102").Value = 49874.75
$ws.Range("J102").Value = 49874.75
$ws.Range("L102").Value = 49874.75
$ws.Range("N102").Value = -56364.75

$ws.Range("H132").Value = 1321.8541
$ws.Range("I132").Value = 1128.3954
$ws.Range("J132").Value = 2985.6
$ws.Range("K132").Value = 3385.1862
$ws.Range("L132").Value = 8956.799999999999
$ws.Range("M132").Value = -855.1862000000001
$ws.Range("N132").Value = -14016.8

$ws.Range("H136").Value = 1353.1
$ws.Range("I136").Value = 1224.1282
$ws.Range("K136").Value = 3672.3846
$ws.Range("M136").Value = -1122.3846
